# Auto-generated Excel COM-interop script
# Applies updated market-price / profit figures to the Leve profit sheets
# (columns H..N) as captured by the scheduled pricing runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2939.4348
$ws.Range("I62").Value = 2312.9167
$ws.Range("K62").Value = 2312.9167
$ws.Range("M62").Value = -1688.9167

$ws.Range("H65").Value = 2939.4348
$ws.Range("I65").Value = 2312.9167
$ws.Range("K65").Value = 11564.5835
$ws.Range("M65").Value = -8444.583500000001

$ws.Range("H129").Value = 748.1892
$ws.Range("I129").Value = 274.5
$ws.Range("J129").Value = 805.6061
$ws.Range("K129").Value = 823.5
$ws.Range("L129").Value = 2416.8183
$ws.Range("M129").Value = 4176.5
$ws.Range("N129").Value = -12416.8183

$ws.Range("H131").Value = 2398.7222
$ws.Range("I131").Value = 1321.3077
$ws.Range("J131").Value = 5200
$ws.Range("K131").Value = 3963.9231
$ws.Range("L131").Value = 15600
$ws.Range("M131").Value = 1076.0769
$ws.Range("N131").Value = -25680

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 410
$ws.Range("J4").Value = 410
$ws.Range("L4").Value = 410
$ws.Range("N4").Value = -642

$ws.Range("H5").Value = 130
$ws.Range("I5").Value = 130
$ws.Range("K5").Value = 130
$ws.Range("M5").Value = -18

$ws.Range("H23").Value = 10000
$ws.Range("I23").Value = 13000
$ws.Range("K23").Value = 13000
$ws.Range("M23").Value = -12741

$ws.Range("H32").Value = 7938.143
$ws.Range("I32").Value = 6144.621
$ws.Range("J32").Value = 18699.273
$ws.Range("K32").Value = 6144.621
$ws.Range("L32").Value = 18699.273
$ws.Range("M32").Value = -5857.621
$ws.Range("N32").Value = -19273.273

$ws.Range("H37").Value = 29850
$ws.Range("J37").Value = 29850
$ws.Range("L37").Value = 29850
$ws.Range("N37").Value = -30396

$ws.Range("H55").Value = 22083.334
$ws.Range("J55").Value = 33055
$ws.Range("L55").Value = 33055
$ws.Range("N55").Value = -33685

$ws.Range("H63").Value = 2843396.2
$ws.Range("I63").Value = 2740
$ws.Range("J63").Value = 15626350
$ws.Range("K63").Value = 2740
$ws.Range("L63").Value = 15626350
$ws.Range("M63").Value = -2054
$ws.Range("N63").Value = -15627722

$ws.Range("H66").Value = 2843396.2
$ws.Range("I66").Value = 2740
$ws.Range("J66").Value = 15626350
$ws.Range("K66").Value = 13700
$ws.Range("L66").Value = 78131750
$ws.Range("M66").Value = -10268
$ws.Range("N66").Value = -78138614

$ws.Range("H80").Value = 51845
$ws.Range("J80").Value = 51845
$ws.Range("L80").Value = 51845
$ws.Range("N80").Value = -53841

$ws.Range("H83").Value = 51845
$ws.Range("J83").Value = 51845
$ws.Range("L83").Value = 155535
$ws.Range("N83").Value = -165519

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 130
$ws.Range("I4").Value = 130
$ws.Range("K4").Value = 130
$ws.Range("M4").Value = -15

$ws.Range("H55").Value = 43780
$ws.Range("J55").Value = 43780
$ws.Range("L55").Value = 43780
$ws.Range("N55").Value = -44326

$ws.Range("H64").Value = 43478980
$ws.Range("I64").Value = 90910120
$ws.Range("J64").Value = 435.75
$ws.Range("K64").Value = 90910120
$ws.Range("L64").Value = 435.75
$ws.Range("M64").Value = -90909895
$ws.Range("N64").Value = -885.75

$ws.Range("H67").Value = 43478980
$ws.Range("I67").Value = 90910120
$ws.Range("J67").Value = 435.75
$ws.Range("K67").Value = 90910120
$ws.Range("L67").Value = 435.75
$ws.Range("M67").Value = -90909340
$ws.Range("N67").Value = -1995.75

$ws.Range("H80").Value = 730.28125
$ws.Range("I80").Value = 769.1667
$ws.Range("K80").Value = 769.1667
$ws.Range("M80").Value = 228.8333

$ws.Range("H83").Value = 730.28125
$ws.Range("I83").Value = 769.1667
$ws.Range("K83").Value = 3845.8335
$ws.Range("M83").Value = 1146.1665

$ws.Range("H111").Value = 40900
$ws.Range("J111").Value = 40900
$ws.Range("L111").Value = 40900
$ws.Range("N111").Value = -49080

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 151.80952
$ws.Range("I22").Value = 149.85715
$ws.Range("J22").Value = 155.71428
$ws.Range("K22").Value = 149.85715
$ws.Range("L22").Value = 155.71428
$ws.Range("M22").Value = 200.14285
$ws.Range("N22").Value = -855.71428

$ws.Range("H31").Value = 4339.744
$ws.Range("I31").Value = 1845.0416
$ws.Range("J31").Value = 7490.9473
$ws.Range("K31").Value = 1845.0416
$ws.Range("L31").Value = 7490.9473
$ws.Range("M31").Value = -1550.0416
$ws.Range("N31").Value = -8080.9473

$ws.Range("H34").Value = 4339.744
$ws.Range("I34").Value = 1845.0416
$ws.Range("J34").Value = 7490.9473
$ws.Range("K34").Value = 1845.0416
$ws.Range("L34").Value = 7490.9473
$ws.Range("M34").Value = -1643.0416
$ws.Range("N34").Value = -7894.9473

$ws.Range("H105").Value = 7813274.5
$ws.Range("J105").Value = 1062
$ws.Range("L105").Value = 1062
$ws.Range("N105").Value = -4556

$ws.Range("H114").Value = 19995
$ws.Range("J114").Value = 19995
$ws.Range("L114").Value = 19995
$ws.Range("N114").Value = -28673

$ws.Range("H122").Value = 1415.0667
$ws.Range("I122").Value = 1093.7
$ws.Range("K122").Value = 3281.1
$ws.Range("M122").Value = -831.1000000000004

$ws.Range("H135").Value = 50780
$ws.Range("J135").Value = 50780
$ws.Range("L135").Value = 50780
$ws.Range("N135").Value = -60920

$ws.Range("H141").Value = 27881.516
$ws.Range("J141").Value = 27881.516
$ws.Range("L141").Value = 27881.516
$ws.Range("N141").Value = -38241.516

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 11188.223
$ws.Range("I2").Value = 16682.334
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 100094.004
$ws.Range("L2").Value = 1200
$ws.Range("M2").Value = -99981.00399999999
$ws.Range("N2").Value = -1426

$ws.Range("H12").Value = 121.5
$ws.Range("J12").Value = 155.33333
$ws.Range("L12").Value = 465.99999
$ws.Range("N12").Value = -811.99999

$ws.Range("H75").Value = 2414.4285
$ws.Range("J75").Value = 2481.3333
$ws.Range("L75").Value = 7443.999899999999
$ws.Range("N75").Value = -9439.999899999999

$ws.Range("H78").Value = 2414.4285
$ws.Range("J78").Value = 2481.3333
$ws.Range("L78").Value = 22331.9997
$ws.Range("N78").Value = -32315.9997

$ws.Range("H131").Value = 721.5599999999999
$ws.Range("J131").Value = 788.4712500000001
$ws.Range("L131").Value = 2365.41375
$ws.Range("N131").Value = -12445.41375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 66.52941
$ws.Range("I2").Value = 72.75
$ws.Range("J2").Value = 51.6
$ws.Range("K2").Value = 72.75
$ws.Range("L2").Value = 51.6
$ws.Range("M2").Value = 40.25
$ws.Range("N2").Value = -277.6

$ws.Range("H97").Value = 1888.579
$ws.Range("I97").Value = 1949
$ws.Range("J97").Value = 1566.3334
$ws.Range("K97").Value = 1949
$ws.Range("L97").Value = 1566.3334
$ws.Range("M97").Value = -1453
$ws.Range("N97").Value = -2558.3334

$ws.Range("H122").Value = 5436.2144
$ws.Range("I122").Value = 5429.5713
$ws.Range("J122").Value = 5442.857
$ws.Range("K122").Value = 16288.7139
$ws.Range("L122").Value = 16328.571
$ws.Range("M122").Value = -13838.7139
$ws.Range("N122").Value = -21228.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H61").Value = 4202.1875
$ws.Range("I61").Value = 1994.1666
$ws.Range("J61").Value = 10826.25
$ws.Range("K61").Value = 1994.1666
$ws.Range("L61").Value = 10826.25
$ws.Range("M61").Value = -1792.1666
$ws.Range("N61").Value = -11230.25

$ws.Range("H93").Value = 1931.9546
$ws.Range("I93").Value = 2043.1428
$ws.Range("J93").Value = 1737.375
$ws.Range("K93").Value = 2043.1428
$ws.Range("L93").Value = 1737.375
$ws.Range("M93").Value = -795.1428000000001
$ws.Range("N93").Value = -4233.375

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws.Range("H113").Value = 4202.1875
$ws.Range("I113").Value = 1994.1666
$ws.Range("J113").Value = 10826.25
$ws.Range("K113").Value = 1994.1666
$ws.Range("L113").Value = 10826.25
$ws.Range("M113").Value = 175.8334
$ws.Range("N113").Value = -15166.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2569.25
$ws.Range("I81").Value = 749.75
$ws.Range("K81").Value = 1499.5
$ws.Range("M81").Value = -438.5

$ws.Range("H84").Value = 2569.25
$ws.Range("I84").Value = 749.75
$ws.Range("K84").Value = 7497.5
$ws.Range("M84").Value = -2193.5

Write-Host "Applied leve-profit updates across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets"
